# Auto-generated edit script applying cryptos.xlsx price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "96.910.69"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +0.42%  "

# Row 3
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "3.691.93"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +0.57%  "

# Row 4
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "236.93"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -1.99%  "

# Row 6
$ws.Range("E6").Value = "  +2.67%  "

# Row 7
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "654.17"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -2.00%  "

# Row 8
$ws.Range("E8").Value = "  -0.23%  "

# Row 9
$ws.Range("E9").Value = "  -0.03%  "

# Row 10
$ws.Range("E10").Value = "  -1.85%  "

# Row 11
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "3.690.17"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +0.57%  "

# Row 12
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "44.15"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -1.73%  "

# Row 13
$ws.Range("E13").Value = "  +2.01%  "

# Row 14
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "0.0000297"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +9.77%  "

# Row 15
$ws.Range("E15").Value = "  +1.20%  "

# Row 16
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "4.379.21"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +0.59%  "

# Row 17
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "96.673.45"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +0.23%  "

# Row 18
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "8.96"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +0.66%  "

# Row 19
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "3.684.81"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +0.59%  "

# Row 20
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "13.07"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +2.62%  "

# Row 21
$ws.Range("E21").Value = "  +1.44%  "

# Row 22
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "0.509"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -4.42%  "

# Row 23
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "522.02"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -0.38%  "

# Row 24
$ws.Range("E24").Value = "  -1.30%  "

# Row 25
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "0.0000212"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +3.60%  "

# Row 26
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "6.91"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -0.43%  "

# Row 27
$ws.Range("B27").Value = "Hedera"
$ws.Range("C27").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "0.202"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +22.35%  "

# Row 28
$ws.Range("B28").Value = "Litecoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "101.49"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -0.97%  "

# Row 29
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "13.35"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +3.04%  "

# Row 30
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "12.37"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +2.34%  "

# Row 31
$ws.Range("E31").Value = "  -1.48%  "

# Row 32
$ws.Range("E32").Value = "  +0.04%  "

# Row 33
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.188"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +1.64%  "

# Row 34
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "1.85"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +1.79%  "

# Row 35
$ws.Range("E35").Value = "  -0.02%  "

# Row 36
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "32.20"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -1.47%  "

# Row 37
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "647.50"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +4.80%  "

# Row 38
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.601"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +2.62%  "

# Row 39
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "8.79"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +0.31%  "

# Row 41
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "6.83"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +10.79%  "

# Row 42
$ws.Range("E42").Value = "  +5.75%  "

# Row 43
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "40.82"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -4.44%  "

# Row 44
$ws.Range("E44").Value = "  +1.10%  "

# Row 45
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.953"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +0.22%  "

# Row 46
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.457"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +7.24%  "

# Row 47
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.0457"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +0.54%  "

# Row 48
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "23.62"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +0.06%  "

# Row 49
$ws.Range("E49").Value = "  -0.74%  "

# Row 50
$ws.Range("E50").Value = "  +0.27%  "

# Row 51
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "3.56"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +0.15%  "
